# Sprint Backlog Burndown - hours edit
# Fills in the burndown "Amount Remaining After..." hours for the tasks
# that were worked on, per the commit message "Edited hours for my tasks
# for the burndown."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("Require compartment to be selected") - Day 1 remaining hours
$ws.Range("D5").Value = 0

# Row 8 ("Add Get Food button functionality") - Day 1..4 remaining hours
$ws.Range("D8").Value = 4.5
$ws.Range("E8").Value = 4.5
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 0

# Reflect the cell the author was last working in when they saved.
$ws.Range("C8").Select()
